$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (England) updates
$ws.Range("F2").Value = 276
$ws.Range("G2").Value = 414
$ws.Range("H2").Value = 317
$ws.Range("I2").Value = 136
$ws.Range("J2").Value = 325
$ws.Range("K2").Value = 589
$ws.Range("L2").Value = 562
$ws.Range("M2").Value = 283
$ws.Range("N2").Value = 287

# Row 3 (Northern Ireland) update
$ws.Range("C3").Value = 107

# Row 4 (Scotland) update
$ws.Range("D4").Value = 622

# Row 5 (Wales) update
$ws.Range("E5").Value = 259
